$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 3.3
$ws.Range("I3").Value = 2.2
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 2
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2.4
$ws.Range("T3").Value = 1.53
$ws.Range("AA3").Value = 2
$ws.Range("AB3").Value = 1.73
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 15
$ws.Range("AE3").Value = 12
$ws.Range("AF3").Value = 34
$ws.Range("AG3").Value = 29
$ws.Range("AI3").Value = 7.5
$ws.Range("AK3").Value = 17
$ws.Range("AM3").Value = 6.5
$ws.Range("AN3").Value = 9.5
$ws.Range("AP3").Value = 21

# Row 4 updates
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.55
$ws.Range("L4").Value = 3.4
$ws.Range("S4").Value = 2.3
$ws.Range("T4").Value = 1.6
$ws.Range("W4").Value = 4.33
$ws.Range("X4").Value = 1.2
$ws.Range("AF4").Value = 29

# Row 5 updates
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 8.5
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 9.5
$ws.Range("AM5").Value = 9.5
